$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-15 (columns A:F) replacing the previous rolling window
$data = @(
    @(45994, 10944.828104324,  10126.2424791726, 17016.26, 6868.28527431307, -0.905510271432225),
    @(45995, 10784.4442086242, 10000.3983874228, 11232.26, 6732.62123910123,  229.198317771833),
    @(45996, 10788.1756069276, 9437.95068976077, 11232.26, 6720.77147257504,  205.269256763992),
    @(45997, 4108.64784199027, 6702.29379470208, 11232.26, 6424.24803408822,  78.9284095329291),
    @(45998, 3743.99018109938, 6433.68825489262, 11232.26, 6161.178280985,    56.7752723282341),
    @(45999, 8386.73995792576, 8683.91322877376, 11232.26, 7357.36552754082,  200.375781513107),
    @(46000, 9706.9828974431,  9620.80077133026, 11232.26, 7769.35234777762,  256.578879962829),
    @(46001, 9706.9828974431,  9395.76998233884, 11232.26, 7769.35234777762,  247.202597088186),
    @(46002, 9706.9828974431,  9285.85435586241, 11232.26, 7769.35234777762,  242.622779318335),
    @(46003, 9706.9828974431,  8607.62748189394, 11232.26, 7769.35234777762,  214.363326236315),
    @(46004, 8488.30621711765, 8555.23921380928, 11232.26, 7366.56059429056,  195.39749200416),
    @(46005, 8386.73995792576, 8608.63650228952, 11232.26, 7357.44211501084,  197.242442387515),
    @(46006, 9746.03784374469, 9418.97915852411, 11232.26, 8042.14436369096,  259.535980092295),
    @(46007, 9746.03784374469, 9661.54570347972, 11232.26, 8042.14436369096,  269.642919465445)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
}
